# Add the new "service_worker-update_available" localization row to the
# KeyValuePairs sheet / "Tabelle2" Excel table (A1:C46 -> A1:C47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# The data lives inside a structured Table ("Tabelle2"); growing it via
# ListRows.Add keeps the table range, AutoFilter range and sheet dimension
# all in sync automatically.
$lo = $ws.ListObjects.Item(1)
$newListRow = $lo.ListRows.Add()
$r = $newListRow.Range.Row

# New Key / English string / German string.
$ws.Cells.Item($r, 1).Value = "service_worker-update_available"
$ws.Cells.Item($r, 2).Value = "A new update is available. Click the button below to refresh the app and get the latest and greatest stuff!"
$ws.Cells.Item($r, 3).Value = "`nEin neues Update ist verfügbar. Klicken Sie auf die Schaltfläche unten, um die App zu aktualisieren und die neuesten und besten Inhalte zu erhalten!"

# Match the existing sheet formatting conventions: column A keeps the plain
# highlighted style used by every other row, while columns B/C (which hold
# the long, wrapped strings) use the wrap-text highlighted style - exactly
# like the other multi-line rows already on the sheet (e.g. rows 14 & 20).
$ws.Range("A45").Copy()
$ws.Range("A" + $r).PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B20:C20").Copy()
$ws.Range("B" + $r + ":C" + $r).PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# Row height grows to fit the wrapped paragraph, matching the other
# multi-line rows' proportions (57pt for this content).
$ws.Rows.Item($r).RowHeight = 57

# Update the window's view to the state left after entering the new row.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
[void]$ws.Range("A" + $r).Select()
